# Refresh crypto price/volume snapshot (GitHub Actions data pull).
# Coin/Link are only rewritten on rows where the ranking reordered;
# Price (D) and Volume(1h) (E) are refreshed for every row below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "66.390.17"
$ws.Range("E2").Value = "  +3.07%  "
# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.185.82"
$ws.Range("E3").Value = "  -0.74%  "
# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.03%  "
# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.82"
$ws.Range("E5").Value = "  +2.04%  "
# Row 6 (Solana)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.44"
$ws.Range("E6").Value = "  -0.27%  "
# Row 7 (XRP)
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
# Row 8 (USDC)
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +0.01%  "
# Row 9 (LidoStakedEther)
$ws.Range("D9").Value = "3.185.41"
$ws.Range("E9").Value = "  -0.57%  "
# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  +2.87%  "
# Row 11 (Toncoin)
$ws.Range("E11").Value = "  -0.68%  "
# Row 12 (Cardano)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.400"
$ws.Range("E12").Value = "  +0.93%  "
# Row 13 (WrappedliquidstakedEther2.0)
$ws.Range("D13").Value = "3.740.21"
$ws.Range("E13").Value = "  -0.53%  "
# Row 14 (TRON)
$ws.Range("E14").Value = "  +1.66%  "
# Row 15 (Avalanche)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.28"
$ws.Range("E15").Value = "  -1.19%  "
# Row 16 (WrappedBTC)
$ws.Range("D16").Value = "66.296.47"
$ws.Range("E16").Value = "  +2.69%  "
# Row 17 (ShibaInu)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000166"
$ws.Range("E17").Value = "  +2.07%  "
# Row 18 (WrappedEther)
$ws.Range("D18").Value = "3.185.95"
$ws.Range("E18").Value = "  -0.62%  "
# Row 19 (Polkadot)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").Value = "  +1.85%  "
# Row 20 (Chainlink)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.06"
$ws.Range("E20").Value = "  +0.46%  "
# Row 21 (BitcoinCash)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.34"
$ws.Range("E21").Value = "  +2.41%  "
# Row 22 (Uniswap)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.35"
$ws.Range("E22").Value = "  +1.89%  "
# Row 23 (Dai)
$ws.Range("E23").Value = "  +0.42%  "
# Row 24 (Litecoin)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.61"
$ws.Range("E24").Value = "  +0.69%  "
# Row 25 (Polygon)
$ws.Range("E25").Value = "  -0.66%  "
# Row 26 (WrappedeETH)
$ws.Range("D26").Value = "3.320.74"
$ws.Range("E26").Value = "  -1.64%  "
# Row 27 (PEPE)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  -1.99%  "
# Row 28 (InternetComputer(DFINITY))
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  +2.98%  "
# Row 29 (Kaspa)
$ws.Range("E29").Value = "  +0.99%  "
# Row 30 (Binance-PegBSC-USD)
$ws.Range("E30").Value = "  +0.18%  "
# Row 31 (USDe)
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.93"
$ws.Range("E31").Value = "  +1.80%  "
# Row 32 (PancakeSwap)
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.49"
$ws.Range("E32").Value = "  -1.67%  "
# Row 33 (NEARProtocol)
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.19%  "
# Row 34 (EthereumClassic)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.18"
$ws.Range("E34").Value = "  -0.15%  "
# Row 35 (Aptos)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.67"
$ws.Range("E35").Value = "  +0.88%  "
# Row 36 (Fetch.AI)
$ws.Range("E36").Value = "  +1.18%  "
# Row 37 (Monero)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "167.15"
$ws.Range("E37").Value = "  +5.30%  "
# Row 38 (ImmutableX)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.48"
$ws.Range("E38").Value = "  +2.38%  "
# Row 39 (Mantle)
$ws.Range("E39").Value = "  +1.70%  "
# Row 40 (Stacks)
$ws.Range("E40").Value = "  +6.88%  "
# Row 41 (EnergySwap)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.28"
$ws.Range("E41").Value = "  +0.00%  "
# Row 42 (dogwifhat)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +2.69%  "
# Row 43 (Maker)
$ws.Range("D43").Value = "2.650.23"
$ws.Range("E43").Value = "  -0.54%  "
# Row 44 (RenderToken)
$ws.Range("E44").Value = "  +2.61%  "
# Row 45 (Filecoin)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.21"
$ws.Range("E45").Value = "  +1.81%  "
# Row 46 (OKB)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.13"
$ws.Range("E46").Value = "  +2.16%  "
# Row 47 (Hedera)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0664"
$ws.Range("E47").Value = "  +2.01%  "
# Row 48 (Bittensor)
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.04"
$ws.Range("E48").Value = "  +2.74%  "
# Row 49 (InjectiveProtocol)
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "326.51"
$ws.Range("E49").Value = "  +1.43%  "
# Row 50 (VeChain)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0276"
$ws.Range("E50").Value = "  +1.91%  "
# Row 51 (Stellar)
$ws.Range("E51").Value = "  +0.75%  "
